# Remove form_id from basic forms
#
# The "settings" sheet has columns: form_title | form_id | version | style | namespaces
# form_id (column B) is being removed entirely, so version/style/namespaces shift
# one column to the left (C->B, D->C, E->D).
#
# Excel's column delete correctly re-indexes cell values/shared-strings, but cell
# *comments* stay anchored to their original absolute addresses, so we manually
# shift the comment text left (re-using AddComment on cells that already own a
# comment keeps the existing "Unknown Author" author instead of minting a new one).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

# Capture the comment text that needs to move left, before anything changes.
$versionComment    = $ws.Range("C1").Comment.Text()
$stylePagesComment = $ws.Range("D1").Comment.Text()
$namespacesComment = $ws.Range("E1").Comment.Text()

# Shift the comments one column left (overwrite in place so the existing
# "Unknown Author" author is reused rather than creating a new one).
$ws.Range("B1").AddComment($versionComment)
$ws.Range("C1").AddComment($stylePagesComment)
$ws.Range("D1").AddComment($namespacesComment)

# The old E1 comment (namespaces) has now been copied onto D1, so drop the
# trailing duplicate.
$ws.Range("E1").Comment.Delete()

# Delete the form_id column itself; this shifts all cell content (and the
# backing shared strings) left automatically.
$ws.Columns("B:B").Delete()

# Restore cursor/selection state to match the post-edit view. Select the
# settings sheet's cell first, then survey's last, so "survey" ends up as
# the active tab again (matches the original tabSelected state).
$ws.Range("A5").Select()

$survey = $wb.Worksheets.Item("survey")
$survey.Range("A12").Select()
